$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.482.52"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.58"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.64"
$ws.Range("E5").Value = "  +5.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4760"
$ws.Range("E7").Value = "  +1.89%  "

$ws.Range("E8").Value = "  +1.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06527"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.88"
$ws.Range("E10").Value = "  +3.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.87"
$ws.Range("E12").Value = "  +3.37%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7395"
$ws.Range("E13").Value = "  +8.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.73"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.128"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "272.77"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.474.14"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.63"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007584"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.127.69"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.229"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.174"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.305"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.34"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.944"
$ws.Range("E28").Value = "  +2.69%  "

$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09978"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.513"
$ws.Range("E31").Value = "  +4.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.318"
$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.061"
$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04779"
$ws.Range("E34").Value = "  +2.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.125"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6994"
$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +1.86%  "

$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.349"
$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.939"
$ws.Range("E41").Value = "  +2.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.08"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4176"
$ws.Range("E43").Value = "  +3.00%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8394"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.72"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.322"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.089"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.51"
$ws.Range("E49").Value = "  +4.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "917.26"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05631"
$ws.Range("E51").Value = "  +0.98%  "
